$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New time-log entry on row 43: clock-out time + the downstream computed
#     columns (Delta Time / Number of minutes / Number of hours / Money).
#     Pull the number formats from row 42 (the previous entry) first, then
#     write the literal clock-out time and the same formulas used by every
#     other row in the log. ---
$ws.Range("C42:G42").Copy()
$ws.Range("C43:G43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C43").Value = 0.76527777777777772
$ws.Range("D43").Formula = "=C43-B43"
$ws.Range("E43").Formula = "=D43*1440"
$ws.Range("F43").Formula = "=E43/60"
$ws.Range("G43").Formula = "=F43*22.5"

# --- Weekly-total bucket (M/N columns) for the week containing the new
#     entry, following the same pattern as the M2:N8 buckets above it. ---
$ws.Range("M8:N8").Copy()
$ws.Range("M9:N9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("M9").Formula = "=SUM(D42:D43)"
$ws.Range("N9").Formula = "=SUM(G42:G43)"

# --- Selection moved to K20 ---
[void]$ws.Range("K20").Select()
